# "Adicionados balanços concatenados em uma única planilha."
# Clears the (redundant/placeholder) numeric 0 values for two rows so the
# cells become blank, matching the already-blank neighboring cells
# (columns B and D) on the same rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C64").ClearContents()
$ws.Range("E64:Z64").ClearContents()

$ws.Range("C79").ClearContents()
$ws.Range("E79:Z79").ClearContents()
